$wb = $excel.ActiveWorkbook

# ALC row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(116, 8).Value = 2666.1667
$ws.Cells.Item(116, 10).Value = 2677.2
$ws.Cells.Item(116, 12).Value = 2677.2
$ws.Cells.Item(116, 14).Value = -9561.200000000001

# ALC row 118
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(118, 8).Value = 1745.6666
$ws.Cells.Item(118, 9).Value = 1745.6666
$ws.Cells.Item(118, 10).Value = 0
$ws.Cells.Item(118, 11).Value = 5236.9998
$ws.Cells.Item(118, 12).Value = 0
$ws.Cells.Item(118, 13).Value = -3579.9998
$ws.Cells.Item(118, 14).ClearContents()

# ALC row 129
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(129, 8).Value = 1079.2826
$ws.Cells.Item(129, 9).Value = 451.88235
$ws.Cells.Item(129, 10).Value = 1447.069
$ws.Cells.Item(129, 11).Value = 1355.64705
$ws.Cells.Item(129, 12).Value = 4341.207
$ws.Cells.Item(129, 13).Value = 3644.35295
$ws.Cells.Item(129, 14).Value = -14341.207

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(138, 8).Value = 1956.5555
$ws.Cells.Item(138, 10).Value = 2323.0293
$ws.Cells.Item(138, 12).Value = 6969.0879
$ws.Cells.Item(138, 14).Value = -17249.0879

# ARM row 82
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(82, 8).Value = 0
$ws.Cells.Item(82, 10).Value = 0
$ws.Cells.Item(82, 12).Value = 0
$ws.Cells.Item(82, 14).ClearContents()

# ARM row 85
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(85, 8).Value = 0
$ws.Cells.Item(85, 10).Value = 0
$ws.Cells.Item(85, 12).Value = 0
$ws.Cells.Item(85, 14).ClearContents()

# BSM row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 1173.2142
$ws.Cells.Item(94, 9).Value = 1111.3636
$ws.Cells.Item(94, 10).Value = 1400
$ws.Cells.Item(94, 11).Value = 1111.3636
$ws.Cells.Item(94, 12).Value = 1400
$ws.Cells.Item(94, 13).Value = -660.3635999999999
$ws.Cells.Item(94, 14).Value = -2302

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 2609.2173
$ws.Cells.Item(134, 9).Value = 2215.1428
$ws.Cells.Item(134, 10).Value = 3222.2222
$ws.Cells.Item(134, 11).Value = 6645.428400000001
$ws.Cells.Item(134, 12).Value = 9666.6666
$ws.Cells.Item(134, 13).Value = -4110.428400000001
$ws.Cells.Item(134, 14).Value = -14736.6666

# CRP row 22
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 572.7273
$ws.Cells.Item(22, 9).Value = 500
$ws.Cells.Item(22, 10).Value = 900
$ws.Cells.Item(22, 11).Value = 500
$ws.Cells.Item(22, 12).Value = 900
$ws.Cells.Item(22, 13).Value = -150
$ws.Cells.Item(22, 14).Value = -1600

# CRP row 44
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(44, 8).Value = 5500
$ws.Cells.Item(44, 9).Value = 5000
$ws.Cells.Item(44, 11).Value = 5000
$ws.Cells.Item(44, 13).Value = -4558

# CRP row 110
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(110, 8).Value = 50000
$ws.Cells.Item(110, 10).Value = 50000
$ws.Cells.Item(110, 12).Value = 50000
$ws.Cells.Item(110, 14).Value = -58180

# CRP row 140
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(140, 8).Value = 39998.57
$ws.Cells.Item(140, 10).Value = 34998.332
$ws.Cells.Item(140, 12).Value = 34998.332
$ws.Cells.Item(140, 14).Value = -45358.332

# CUL row 3
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(3, 8).Value = 5750
$ws.Cells.Item(3, 9).Value = 5750
$ws.Cells.Item(3, 11).Value = 17250
$ws.Cells.Item(3, 13).Value = -17138

# CUL row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 611.6818
$ws.Cells.Item(5, 9).Value = 430.1579
$ws.Cells.Item(5, 10).Value = 1761.3334
$ws.Cells.Item(5, 11).Value = 1290.4737
$ws.Cells.Item(5, 12).Value = 5284.0002
$ws.Cells.Item(5, 13).Value = -1178.4737
$ws.Cells.Item(5, 14).Value = -5508.0002

# CUL row 24
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(24, 8).Value = 1200
$ws.Cells.Item(24, 9).Value = 0
$ws.Cells.Item(24, 10).Value = 1200
$ws.Cells.Item(24, 11).Value = 0
$ws.Cells.Item(24, 12).Value = 3600
$ws.Cells.Item(24, 13).ClearContents()
$ws.Cells.Item(24, 14).Value = -4060

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 1031.5
$ws.Cells.Item(131, 10).Value = 1142.1482
$ws.Cells.Item(131, 12).Value = 3426.4446
$ws.Cells.Item(131, 14).Value = -13506.4446

# CUL row 132
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(132, 8).Value = 1665.32
$ws.Cells.Item(132, 9).Value = 1184.3334
$ws.Cells.Item(132, 10).Value = 1935.875
$ws.Cells.Item(132, 11).Value = 10659.0006
$ws.Cells.Item(132, 12).Value = 17422.875
$ws.Cells.Item(132, 13).Value = -8129.000599999999
$ws.Cells.Item(132, 14).Value = -22482.875

# CUL row 133
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(133, 8).Value = 14169.692
$ws.Cells.Item(133, 9).Value = 7632
$ws.Cells.Item(133, 10).Value = 18255.75
$ws.Cells.Item(133, 11).Value = 22896
$ws.Cells.Item(133, 12).Value = 54767.25
$ws.Cells.Item(133, 13).Value = -17836
$ws.Cells.Item(133, 14).Value = -64887.25

# CUL row 135
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(135, 8).Value = 611.6818
$ws.Cells.Item(135, 9).Value = 430.1579
$ws.Cells.Item(135, 10).Value = 1761.3334
$ws.Cells.Item(135, 11).Value = 3871.4211
$ws.Cells.Item(135, 12).Value = 15852.0006
$ws.Cells.Item(135, 13).Value = -1336.4211
$ws.Cells.Item(135, 14).Value = -20922.0006

# GSM row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 5529.091
$ws.Cells.Item(70, 9).Value = 5585.75
$ws.Cells.Item(70, 10).Value = 5429.9375
$ws.Cells.Item(70, 11).Value = 5585.75
$ws.Cells.Item(70, 12).Value = 5429.9375
$ws.Cells.Item(70, 13).Value = -5315.75
$ws.Cells.Item(70, 14).Value = -5969.9375

# GSM row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(73, 8).Value = 5529.091
$ws.Cells.Item(73, 9).Value = 5585.75
$ws.Cells.Item(73, 10).Value = 5429.9375
$ws.Cells.Item(73, 11).Value = 5585.75
$ws.Cells.Item(73, 12).Value = 5429.9375
$ws.Cells.Item(73, 13).Value = -4649.75
$ws.Cells.Item(73, 14).Value = -7301.9375

# GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 1513.6
$ws.Cells.Item(122, 9).Value = 1563.6666
$ws.Cells.Item(122, 10).Value = 1250.75
$ws.Cells.Item(122, 11).Value = 4690.9998
$ws.Cells.Item(122, 12).Value = 3752.25
$ws.Cells.Item(122, 13).Value = -2240.9998
$ws.Cells.Item(122, 14).Value = -8652.25

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 2744.5757
$ws.Cells.Item(132, 9).Value = 2451.389
$ws.Cells.Item(132, 10).Value = 3096.4
$ws.Cells.Item(132, 11).Value = 7354.167
$ws.Cells.Item(132, 12).Value = 9289.200000000001
$ws.Cells.Item(132, 13).Value = -4824.167
$ws.Cells.Item(132, 14).Value = -14349.2

# LTW row 16
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 11906807
$ws.Cells.Item(16, 9).Value = 1854.6
$ws.Cells.Item(16, 10).Value = 71431570
$ws.Cells.Item(16, 11).Value = 1854.6
$ws.Cells.Item(16, 12).Value = 71431570
$ws.Cells.Item(16, 13).Value = -1684.6
$ws.Cells.Item(16, 14).Value = -71431910

# LTW row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 5871
$ws.Cells.Item(22, 9).Value = 1233.3334
$ws.Cells.Item(22, 10).Value = 8479.6875
$ws.Cells.Item(22, 11).Value = 1233.3334
$ws.Cells.Item(22, 12).Value = 8479.6875
$ws.Cells.Item(22, 13).Value = -938.3334
$ws.Cells.Item(22, 14).Value = -9069.6875

# LTW row 27
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(27, 8).Value = 5871
$ws.Cells.Item(27, 9).Value = 1233.3334
$ws.Cells.Item(27, 10).Value = 8479.6875
$ws.Cells.Item(27, 11).Value = 1233.3334
$ws.Cells.Item(27, 12).Value = 8479.6875
$ws.Cells.Item(27, 13).Value = -1126.3334
$ws.Cells.Item(27, 14).Value = -8693.6875

# WVR row 47
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(47, 8).Value = 0
$ws.Cells.Item(47, 10).Value = 0
$ws.Cells.Item(47, 12).Value = 0
$ws.Cells.Item(47, 14).ClearContents()

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 2865.0334
$ws.Cells.Item(136, 9).Value = 2480.611
$ws.Cells.Item(136, 10).Value = 3441.6667
$ws.Cells.Item(136, 11).Value = 7441.833
$ws.Cells.Item(136, 12).Value = 10325.0001
$ws.Cells.Item(136, 13).Value = -4891.833
$ws.Cells.Item(136, 14).Value = -15425.0001
